$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.74184335344219499
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0.16500000000000001
$ws.Range("B3").Value = 0.69910893998062895
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.154
$ws.Range("B4").Value = 0.66521175214108597
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0.14299999999999999
$ws.Range("B5").Value = 0.65216376659028397
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.13600000000000001
$ws.Range("B6").Value = 0.64511095935346396
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0.123
$ws.Range("B7").Value = 0.65371835032927095
$ws.Range("C7").Value = 0.026700000000000002
$ws.Range("D7").Value = 0.096000000000000002
$ws.Range("B8").Value = 0.696800520862071
$ws.Range("C8").Value = 0.083150000000000002
$ws.Range("D8").Value = 0.068000000000000005
$ws.Range("B9").Value = 0.75299294933528904
$ws.Range("C9").Value = 0.27826000000000001
$ws.Range("D9").Value = 0.045999999999999999
$ws.Range("B10").Value = 0.81319813491132298
$ws.Range("C10").Value = 0.47961000000000004
$ws.Range("D10").Value = 0.047
$ws.Range("B11").Value = 0.86911679303537503
$ws.Range("C11").Value = 0.64805999999999997
$ws.Range("D11").Value = 0.1
$ws.Range("B12").Value = 0.90369221568139502
$ws.Range("C12").Value = 0.76993
$ws.Range("D12").Value = 0.14699999999999999
$ws.Range("B13").Value = 0.91808220183594402
$ws.Range("C13").Value = 0.83309
$ws.Range("D13").Value = 0.20799999999999999
$ws.Range("B14").Value = 0.93924668940026201
$ws.Range("C14").Value = 0.84025000000000005
$ws.Range("D14").Value = 0.29299999999999998
$ws.Range("B15").Value = 0.94932964106478401
$ws.Range("C15").Value = 0.79170000000000007
$ws.Range("D15").Value = 0.38600000000000001
$ws.Range("B16").Value = 0.92978392955529099
$ws.Range("C16").Value = 0.68528999999999995
$ws.Range("D16").Value = 0.46500000000000002
$ws.Range("B17").Value = 0.90963753156365601
$ws.Range("C17").Value = 0.53488000000000002
$ws.Range("D17").Value = 0.52200000000000002
$ws.Range("B18").Value = 0.90548461984667705
$ws.Range("C18").Value = 0.33491000000000004
$ws.Range("D18").Value = 0.57299999999999995
$ws.Range("B19").Value = 0.903493266804829
$ws.Range("C19").Value = 0.12425
$ws.Range("D19").Value = 0.61499999999999999
$ws.Range("B20").Value = 0.88605727596248696
$ws.Range("C20").Value = 0.038609999999999998
$ws.Range("D20").Value = 0.64300000000000002
$ws.Range("B21").Value = 0.87104631888789097
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0.63700000000000001
$ws.Range("B22").Value = 0.84917763577558203
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0.57399999999999995
$ws.Range("B23").Value = 0.82627007780236605
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0.497
$ws.Range("B24").Value = 0.80277850053261202
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0.432
$ws.Range("B25").Value = 0.75405035565798495
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 0.371

$ws.Range("E2:H26").Select()
